$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Network1")

# Update D12:D20 values
$ws.Range("D12").Value = 13
$ws.Range("D13").Value = 15
$ws.Range("D14").Value = 17
$ws.Range("D15").Value = 19
$ws.Range("D16").Value = 26
$ws.Range("D17").Value = 37
$ws.Range("D18").Value = 38
$ws.Range("D19").Value = 49
$ws.Range("D20").Value = 50

# Update selection to D21
$ws.Activate()
$ws.Range("D21").Select()
